$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns C and D for rows 1-5
for ($row = 1; $row -le 5; $row++) {
    $cCell = $ws.Cells.Item($row, 3)  # Column C
    $dCell = $ws.Cells.Item($row, 4)  # Column D

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
